# "6 month coverage effective from 2026 for scenario 2"
#
# Both worksheets currently lay out one calendar year per column (every
# year from 2018 to 2040). The edit changes the timeline to 6-month
# (half-year) steps over the same span, and extends/refills the
# coverage (sheet1) / market-share (sheet2) rows to match the new,
# wider column grid.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Platform Coverage"
$ws2 = $wb.Worksheets.Item(2)   # "MarketShare"

# ---------------------------------------------------------------
# Sheet 1 : "Platform Coverage"
# Year header used to start at column H (col 8) = 2018, one column
# per year, through column AD (col 30) = 2040.
# It now starts at column H (col 8) = 2018, one column per HALF year,
# through column AZ (col 52) = 2040.
# ---------------------------------------------------------------

# Row 1: rewrite the whole header run with half-year steps.
$col = 8
$year = 2018.0
while ($year -le 2040.0) {
    $ws1.Cells.Item(1, $col).Value = $year
    $col = $col + 1
    $year = $year + 0.5
}

# Row 2 (age 5-15, product routine coverage 0.6): previously had values
# every other column from 2018 to 2024 (H,J,L,N). Those same letters
# still mean 2018/2019/2020/2021 under the new half-year grid, so they
# are left untouched; we extend the same 0.6 coverage through the next
# whole-year columns (2022-2025 => P,R,T,V).
foreach ($c in @("P", "R", "T", "V")) {
    $ws1.Range($c + "2").Value = 0.6
}

# Rows 3-5 (age 2-15 / 15-50 / 50-65): previously had coverage every
# other year from 2026 to 2040 (P,R,T,V,X,Z,AB,AD in the old grid).
# Under the new half-year grid those letters no longer line up with
# the same years, so the stale cells left of column X are cleared, and
# the row is refilled with the SAME value every half-year column from
# X (2026) through AZ (2040) -- i.e. 6-month coverage starting 2026.
foreach ($c in @("P", "R", "T", "V")) {
    $ws1.Range($c + "3").ClearContents() | Out-Null
    $ws1.Range($c + "4").ClearContents() | Out-Null
    $ws1.Range($c + "5").ClearContents() | Out-Null
}

$rowValues = @{ 3 = 0.8; 4 = 0.5; 5 = 0.5 }
foreach ($r in $rowValues.Keys) {
    $v = $rowValues[$r]
    for ($col = 24; $col -le 52; $col++) {
        # columns X (24) .. AZ (52)
        $ws1.Cells.Item($r, $col).Value = $v
    }
}

# ---------------------------------------------------------------
# Sheet 2 : "MarketShare"
# Year header used to start at column D (col 4) = 2018, one column
# per year, through column Z (col 26) = 2040.
# It now starts at column D (col 4) = 2018, one column per HALF year,
# through column AV (col 48) = 2040.
# ---------------------------------------------------------------

# Row 1: rewrite the whole header run with half-year steps.
$col = 4
$year = 2018.0
while ($year -le 2040.0) {
    $ws2.Cells.Item(1, $col).Value = $year
    $col = $col + 1
    $year = $year + 0.5
}

# Row 3: was 1 in every column from D (2018) to Z (2040, old 1yr grid).
# Those columns keep the same meaning under the new grid, so just
# extend the same value of 1 through the newly added half-year columns
# (AA .. AV).
for ($col = 27; $col -le 48; $col++) {
    # columns AA (27) .. AV (48)
    $ws2.Cells.Item(3, $col).Value = 1
}

# ---------------------------------------------------------------
# View state: restore each sheet's selection the way it was left in
# the saved file, and make sure sheet 1 ends up the active/selected
# tab again (matches the workbook's tabSelected state).
# ---------------------------------------------------------------

$ws2.Activate() | Out-Null
$ws2.Range("D1:AV1").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("AW3:AW5").Select() | Out-Null
